$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1: move selection to A4 (tabSelected will move off automatically once Sheet2 becomes active)
$ws1.Range("A4").Select() | Out-Null

# Add Sheet2 right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Populate header row
$ws2.Range("A1").Value = "Make"
$ws2.Range("B1").Value = "Model"
$ws2.Range("C1").Value = "Type"

# Populate data row - set B2 before A2 so shared-string order matches (v70 before volvo)
$ws2.Range("B2").Value = "v70"
$ws2.Range("A2").Value = "volvo"
$ws2.Range("C2").Value = "1fb90eb32df417632de158bfe4c24089"

# Style C2 with the Consolas font used for the review hash
$ws2.Range("C2").Font.Name = "Consolas"
$ws2.Range("C2").Font.Size = 9
$ws2.Range("C2").Font.Color = 2236962

# Column C is widened to fit the hash text
$ws2.Columns.Item(3).ColumnWidth = 33.65

# Page setup matching the source workbook
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Select C2 and make Sheet2 the active sheet/tab
$ws2.Range("C2").Select() | Out-Null
